# Auto-generated edit script: applies the scheduled-runner updates
# to the Rafflesia_Profits workbook's per-item profit calculations.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 244
$ws.Range("J12").Value = 275
$ws.Range("L12").Value = 275
$ws.Range("N12").Value = -615
$ws.Range("H28").Value = 808.1
$ws.Range("I28").Value = 404.57144
$ws.Range("K28").Value = 404.57144
$ws.Range("M28").Value = 80.42856
$ws.Range("H48").Value = 6973
$ws.Range("J48").Value = 6973
$ws.Range("L48").Value = 20919
$ws.Range("N48").Value = -21503
$ws.Range("H56").Value = 6973
$ws.Range("J56").Value = 6973
$ws.Range("L56").Value = 20919
$ws.Range("N56").Value = -21987
$ws.Range("H58").Value = 1299
$ws.Range("J58").Value = 1223.75
$ws.Range("L58").Value = 3671.25
$ws.Range("N58").Value = -3971.25
$ws.Range("H132").Value = 7976.6665
$ws.Range("I132").Value = 7976.6665
$ws.Range("K132").Value = 23929.9995
$ws.Range("M132").Value = -21399.9995
$ws.Range("H137").Value = 3415.5557
$ws.Range("I137").Value = 3467.5
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 10402.5
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -7852.5
$ws.Range("N137").Value = -14100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3402.8125
$ws.Range("I32").Value = 3353.2144
$ws.Range("J32").Value = 3750
$ws.Range("K32").Value = 3353.2144
$ws.Range("L32").Value = 3750
$ws.Range("M32").Value = -3066.2144
$ws.Range("N32").Value = -4324
$ws.Range("H61").Value = 450
$ws.Range("I61").Value = 450
$ws.Range("K61").Value = 450
$ws.Range("M61").Value = -238
$ws.Range("H119").Value = 49500.5
$ws.Range("J119").Value = 49500.5
$ws.Range("L119").Value = 49500.5
$ws.Range("N119").Value = -59176.5
$ws.Range("H124").Value = 40000.5
$ws.Range("J124").Value = 40000.5
$ws.Range("L124").Value = 40000.5
$ws.Range("N124").Value = -49820.5
$ws.Range("H132").Value = 1800
$ws.Range("I132").Value = 1800
$ws.Range("K132").Value = 5400
$ws.Range("M132").Value = -2870
$ws.Range("H136").Value = 450
$ws.Range("I136").Value = 450
$ws.Range("K136").Value = 1350
$ws.Range("M136").Value = 1200
$ws.Range("H139").Value = 80238.336
$ws.Range("J139").Value = 80238.336
$ws.Range("L139").Value = 80238.336
$ws.Range("N139").Value = -90518.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 29999
$ws.Range("J21").Value = 29999
$ws.Range("L21").Value = 29999
$ws.Range("N21").Value = -30471
$ws.Range("H94").Value = 519.8
$ws.Range("I94").Value = 437.25
$ws.Range("K94").Value = 437.25
$ws.Range("M94").Value = 13.75
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H37").Value = 4550
$ws.Range("I37").Value = 4550
$ws.Range("K37").Value = 4550
$ws.Range("M37").Value = -4443
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H45").Value = 5500
$ws.Range("I45").Value = 5500
$ws.Range("K45").Value = 5500
$ws.Range("M45").Value = -4907
$ws.Range("H132").Value = 7726.1816
$ws.Range("J132").Value = 9571.429
$ws.Range("L132").Value = 28714.287
$ws.Range("N132").Value = -33774.287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 683.5714
$ws.Range("I5").Value = 511.5
$ws.Range("J5").Value = 752.4
$ws.Range("K5").Value = 1534.5
$ws.Range("L5").Value = 2257.2
$ws.Range("M5").Value = -1422.5
$ws.Range("N5").Value = -2481.2
$ws.Range("H45").Value = 2030
$ws.Range("I45").Value = 2030
$ws.Range("K45").Value = 6090
$ws.Range("M45").Value = -5558
$ws.Range("H60").Value = 246.5
$ws.Range("I60").Value = 196
$ws.Range("J60").Value = 499
$ws.Range("K60").Value = 588
$ws.Range("L60").Value = 1497
$ws.Range("M60").Value = -337
$ws.Range("N60").Value = -1999
$ws.Range("H68").Value = 1166.6666
$ws.Range("J68").Value = 2500
$ws.Range("L68").Value = 7500
$ws.Range("N68").Value = -9122
$ws.Range("H71").Value = 1166.6666
$ws.Range("J71").Value = 2500
$ws.Range("L71").Value = 22500
$ws.Range("N71").Value = -30612
$ws.Range("H135").Value = 683.5714
$ws.Range("I135").Value = 511.5
$ws.Range("J135").Value = 752.4
$ws.Range("K135").Value = 4603.5
$ws.Range("L135").Value = 6771.599999999999
$ws.Range("M135").Value = -2068.5
$ws.Range("N135").Value = -11841.6
$ws.Range("H138").Value = 1666.6666
$ws.Range("J138").Value = 2000
$ws.Range("L138").Value = 6000
$ws.Range("N138").Value = -16280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H9").Value = 3336666.8
$ws.Range("J9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("N9").Value = -5340
$ws.Range("H11").Value = 4736066.5
$ws.Range("I11").Value = 7224666.5
$ws.Range("K11").Value = 7224666.5
$ws.Range("M11").Value = -7224527.5
$ws.Range("H43").Value = 2898.25
$ws.Range("I43").Value = 2699.6667
$ws.Range("J43").Value = 3017.4
$ws.Range("K43").Value = 2699.6667
$ws.Range("L43").Value = 3017.4
$ws.Range("M43").Value = -2548.6667
$ws.Range("N43").Value = -3319.4
$ws.Range("H118").Value = 40000
$ws.Range("J118").Value = 40000
$ws.Range("L118").Value = 40000
$ws.Range("N118").Value = -43314
$ws.Range("H122").Value = 1516.5
$ws.Range("I122").Value = 1516.5
$ws.Range("K122").Value = 4549.5
$ws.Range("M122").Value = -2099.5
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 420
$ws.Range("H27").Value = 420
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4071.4285
$ws.Range("I2").Value = 4416.6665
$ws.Range("K2").Value = 4416.6665
$ws.Range("M2").Value = -4304.6665
$ws.Range("H52").Value = 30021
$ws.Range("I52").Value = 30021
$ws.Range("K52").Value = 30021
$ws.Range("M52").Value = -29795
$ws.Range("H113").Value = 567.6667
$ws.Range("I113").Value = 363.375
$ws.Range("J113").Value = 1221.4
$ws.Range("K113").Value = 1090.125
$ws.Range("L113").Value = 3664.2
$ws.Range("M113").Value = 1079.875
$ws.Range("N113").Value = -8004.200000000001
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H132").Value = 4337.6
$ws.Range("I132").Value = 4337.6
$ws.Range("K132").Value = 13012.8
$ws.Range("M132").Value = -10482.8

Write-Output "Applied scheduled profit updates to ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR."